$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# New rows to append to Table1 (Feature, ID, Requirment)
# ---------------------------------------------------------------------------
$newRowsData = @(
    @("Notifications",    "LH_CRS_NOTIFICATIONS_004",       "users can follow specific categories to receive notifications when new content is added"),
    @("System Constrain", "LH_CRS_SYSTEM-CONSTRAINS_005",   "web-based system/PC based "),
    @("ID Constrain",     "LH_CRS_ID-CONSTRAINS_006",       "each user must be assigned a unique user ID after creating an account"),
    @("ADMIN Constrain",  "LH_CRS_ADMIN-CONSTRAINS_007",    "admin features must be available for content and user management")
)

foreach ($rowVals in $newRowsData) {
    $lastRow = $tbl.ListRows.Item($tbl.ListRows.Count)
    $lastRow.Range.Copy()

    $newRow = $tbl.ListRows.Add()
    $r = $newRow.Range
    $r.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $r.Item(1, 1).Value = $rowVals[0]
    $r.Item(1, 2).Value = $rowVals[1]
    $r.Item(1, 3).Value = $rowVals[2]
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row-level formatting tweaks that mirror the tracked-changes source file
# ---------------------------------------------------------------------------

# "System Constrain" row -> Feature cell carries an explicit black font colour
$ws.Range("A6").Font.Color = 0

# "ID Constrain" row -> Requirment cell lost its table border formatting
$ws.Range("C7").ClearFormats()
$ws.Range("C7").Value = "each user must be assigned a unique user ID after creating an account"

# "ADMIN Constrain" row is the new last row of the table, so it only keeps a
# top border (no bottom) on each side, and its ID cell renders bold.
$adminFeature = $ws.Range("A8")
$adminId = $ws.Range("B8")
$adminReq = $ws.Range("C8")

$adminFeature.Borders.Item(9).LineStyle = 0    # xlEdgeBottom -> none

$adminId.Font.Bold = $true
$adminId.Borders.Item(7).LineStyle = 0         # xlEdgeLeft -> none
$adminId.Borders.Item(9).LineStyle = 0         # xlEdgeBottom -> none

$adminReq.Borders.Item(10).LineStyle = 0       # xlEdgeRight -> none
$adminReq.Borders.Item(9).LineStyle = 0        # xlEdgeBottom -> none

# ---------------------------------------------------------------------------
# Data validation (blank-allowed, in-cell dropdown w/ no source list) carries
# over to the new "Notifications"/"ID Constrain"/"ADMIN Constrain" rows
# ---------------------------------------------------------------------------
foreach ($addr in @("A5", "A7:A8")) {
    $rng = $ws.Range($addr)
    $v = $rng.Validation
    $v.Add(0) | Out-Null
    $v.IgnoreBlank = $true
    $v.InCellDropdown = $false
    $v.ShowInput = $false
    $v.ShowError = $false
}

# ---------------------------------------------------------------------------
# Misc cosmetic bits from the diff
# ---------------------------------------------------------------------------
$ws.Range("B12").Select()
